# Updated cryptos list values (Price and Volume(1h) columns)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '23.445.97'
$ws.Range("E2").Value = '  +1.91%  '

$ws.Range("D3").Value = '1.628.72'
$ws.Range("E3").Value = '  +2.79%  '

$ws.Range("D4").Value = "'0.9964"
$ws.Range("E4").Value = '  -0.75%  '

$ws.Range("D5").Value = "'307.39"
$ws.Range("E5").Value = '  +2.28%  '

$ws.Range("D6").Value = "'0.9968"
$ws.Range("E6").Value = '  -0.64%  '

$ws.Range("D7").Value = "'0.3782"
$ws.Range("E7").Value = '  +0.58%  '

$ws.Range("D8").Value = "'53.12"
$ws.Range("E8").Value = '  +5.16%  '

$ws.Range("D9").Value = "'0.3667"
$ws.Range("E9").Value = '  +2.54%  '

$ws.Range("D10").Value = "'1.285"
$ws.Range("E10").Value = '  +5.80%  '

$ws.Range("D11").Value = "'0.08203"
$ws.Range("E11").Value = '  +2.37%  '

$ws.Range("D12").Value = "'0.9963"
$ws.Range("E12").Value = '  -0.76%  '

$ws.Range("D13").Value = "'23.31"
$ws.Range("E13").Value = '  +6.58%  '

$ws.Range("D14").Value = "'6.678"
$ws.Range("E14").Value = '  +3.50%  '

$ws.Range("D15").Value = "'7.468"
$ws.Range("E15").Value = '  +2.37%  '

$ws.Range("D16").Value = "'0.00001266"
$ws.Range("E16").Value = '  +3.86%  '

$ws.Range("D17").Value = '1.630.55'
$ws.Range("E17").Value = '  +2.66%  '

$ws.Range("D18").Value = "'94.85"
$ws.Range("E18").Value = '  +2.96%  '

$ws.Range("D19").Value = "'0.06939"
$ws.Range("E19").Value = '  +2.18%  '

$ws.Range("D20").Value = "'18.43"
$ws.Range("E20").Value = '  +3.06%  '

$ws.Range("D21").Value = "'6.591"
$ws.Range("E21").Value = '  +2.68%  '

$ws.Range("D22").Value = "'0.9991"
$ws.Range("E22").Value = '  -0.41%  '

$ws.Range("D23").Value = "'13.01"
$ws.Range("E23").Value = '  +1.77%  '

$ws.Range("D24").Value = '23.470.14'
$ws.Range("E24").Value = '  +2.01%  '

$ws.Range("D25").Value = "'3.135"
$ws.Range("E25").Value = '  +13.78%  '

$ws.Range("D26").Value = "'2.432"
$ws.Range("E26").Value = '  +2.82%  '

$ws.Range("D27").Value = "'21.41"

$ws.Range("D28").Value = "'150.69"

$ws.Range("D29").Value = "'5.292"
$ws.Range("E29").Value = '  +1.81%  '

$ws.Range("D30").Value = "'136.39"
$ws.Range("E30").Value = '  +2.88%  '

$ws.Range("D31").Value = "'2.423"
$ws.Range("E31").Value = '  +3.34%  '

$ws.Range("D32").Value = "'6.959"
$ws.Range("E32").Value = '  +7.16%  '

$ws.Range("D33").Value = '1.806.67'
$ws.Range("E33").Value = '  +2.37%  '

$ws.Range("D34").Value = "'0.9783"
$ws.Range("E34").Value = '  +4.46%  '

$ws.Range("D35").Value = "'0.02807"
$ws.Range("E35").Value = '  +5.07%  '

$ws.Range("D36").Value = "'10.50"
$ws.Range("E36").Value = '  +5.15%  '

$ws.Range("D37").Value = "'0.07496"
$ws.Range("E37").Value = '  +2.19%  '

$ws.Range("D38").Value = "'6.253"
$ws.Range("E38").Value = '  +3.34%  '

$ws.Range("D39").Value = "'0.2539"
$ws.Range("E39").Value = '  +2.77%  '

$ws.Range("D40").Value = "'0.08844"
$ws.Range("E40").Value = '  +1.07%  '

$ws.Range("D41").Value = "'1.409"
$ws.Range("E41").Value = '  +5.57%  '

$ws.Range("D42").Value = "'0.7177"
$ws.Range("E42").Value = '  +4.53%  '

$ws.Range("D43").Value = "'12.82"
$ws.Range("E43").Value = '  +8.19%  '

$ws.Range("D44").Value = "'16.20"
$ws.Range("E44").Value = '  +8.62%  '

$ws.Range("D45").Value = "'0.6633"
$ws.Range("E45").Value = '  +3.89%  '

$ws.Range("E46").Value = '  +5.49%  '

$ws.Range("D47").Value = "'4.036"
$ws.Range("E47").Value = '  +1.31%  '

$ws.Range("E48").Value = '  -0.64%  '

$ws.Range("D49").Value = "'0.08028"
$ws.Range("E49").Value = '  +1.94%  '

$ws.Range("D50").Value = "'132.06"
$ws.Range("E50").Value = '  +1.10%  '

$ws.Range("E51").Value = '  +2.28%  '
